$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plants")

# Mark the "TMin" variable (row 14) as not to be tuned.
$ws.Range("E14").Value = "no"

# Rename the header columns H1:K1 to reflect the proper "Min"/"Max" naming
# scheme instead of the old "Border" naming scheme, and drop the bold/wrap
# header style those cells carried previously (style "s=1" removed in diff).
$ws.Range("H1:K1").Style = "Normal"
$ws.Range("H1").Value = "Hard Min"
$ws.Range("I1").Value = "Soft Min"
$ws.Range("J1").Value = "Soft Max"
$ws.Range("K1").Value = "Hard Max"
